# Scheduled-runner data refresh: updates cached market-board figures
# (currentAveragePrice*, Leve price/profit columns) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 281.06668
$ws.Range("I33").Value = 236.85715
$ws.Range("K33").Value = 236.85715
$ws.Range("M33").Value = -7.85714999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 32903.945
$ws.Range("I51").Value = 6985
$ws.Range("J51").Value = 34428.59
$ws.Range("K51").Value = 6985
$ws.Range("L51").Value = 34428.59
$ws.Range("M51").Value = -6501
$ws.Range("N51").Value = -35396.59

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 3005.6667
$ws.Range("J58").Value = 5714.143
$ws.Range("L58").Value = 17142.429
$ws.Range("N58").Value = -17442.429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 8821792
$ws.Range("I106").Value = 9498469
$ws.Range("K106").Value = 9498469
$ws.Range("M106").Value = -9497838

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3026.111
$ws.Range("I125").Value = 609.75
$ws.Range("K125").Value = 5487.75
$ws.Range("M125").Value = -3027.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3901.9
$ws.Range("I131").Value = 2002.7142
$ws.Range("K131").Value = 6008.142599999999
$ws.Range("M131").Value = -968.1425999999992

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5644035
$ws.Range("J137").Value = 9625745
$ws.Range("L137").Value = 28877235
$ws.Range("N137").Value = -28882335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5959.7295
$ws.Range("I138").Value = 2026.5555
$ws.Range("J138").Value = 7223.9644
$ws.Range("K138").Value = 6079.666499999999
$ws.Range("L138").Value = 21671.8932
$ws.Range("M138").Value = -939.6664999999994
$ws.Range("N138").Value = -31951.8932

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 94793.336
$ws.Range("I2").Value = 17856.285
$ws.Range("J2").Value = 202505.2
$ws.Range("K2").Value = 17856.285
$ws.Range("L2").Value = 202505.2
$ws.Range("M2").Value = -17743.285
$ws.Range("N2").Value = -202731.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1783.5834
$ws.Range("I32").Value = 1708.6868
$ws.Range("K32").Value = 1708.6868
$ws.Range("M32").Value = -1421.6868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5403.7896
$ws.Range("I61").Value = 5502.617
$ws.Range("K61").Value = 5502.617
$ws.Range("M61").Value = -5290.617

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9833012
$ws.Range("I74").Value = 12230923
$ws.Range("J74").Value = 1577.8
$ws.Range("K74").Value = 12230923
$ws.Range("L74").Value = 1577.8
$ws.Range("M74").Value = -12230049
$ws.Range("N74").Value = -3325.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9833012
$ws.Range("I77").Value = 12230923
$ws.Range("J77").Value = 1577.8
$ws.Range("K77").Value = 61154615
$ws.Range("L77").Value = 7889
$ws.Range("M77").Value = -61150247
$ws.Range("N77").Value = -16625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 94793.336
$ws.Range("I116").Value = 17856.285
$ws.Range("J116").Value = 202505.2
$ws.Range("K116").Value = 17856.285
$ws.Range("L116").Value = 202505.2
$ws.Range("M116").Value = -15562.285
$ws.Range("N116").Value = -207093.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 362522.2
$ws.Range("I122").Value = 2633.3447
$ws.Range("K122").Value = 7900.034100000001
$ws.Range("M122").Value = -5450.034100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3499.4856
$ws.Range("I132").Value = 3017.3928
$ws.Range("K132").Value = 9052.178400000001
$ws.Range("M132").Value = -6522.178400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5403.7896
$ws.Range("I136").Value = 5502.617
$ws.Range("K136").Value = 16507.851
$ws.Range("M136").Value = -13957.851

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 94793.336
$ws.Range("I3").Value = 17856.285
$ws.Range("J3").Value = 202505.2
$ws.Range("K3").Value = 17856.285
$ws.Range("L3").Value = 202505.2
$ws.Range("M3").Value = -17742.285
$ws.Range("N3").Value = -202733.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 593.4706
$ws.Range("I80").Value = 255.8
$ws.Range("J80").Value = 734.1667
$ws.Range("K80").Value = 255.8
$ws.Range("L80").Value = 734.1667
$ws.Range("M80").Value = 742.2
$ws.Range("N80").Value = -2730.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 593.4706
$ws.Range("I83").Value = 255.8
$ws.Range("J83").Value = 734.1667
$ws.Range("K83").Value = 1279
$ws.Range("L83").Value = 3670.8335
$ws.Range("M83").Value = 3713
$ws.Range("N83").Value = -13654.8335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7018.4443
$ws.Range("I86").Value = 8345.23
$ws.Range("K86").Value = 8345.23
$ws.Range("M86").Value = -7222.23

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 7018.4443
$ws.Range("I89").Value = 8345.23
$ws.Range("K89").Value = 41726.14999999999
$ws.Range("M89").Value = -36110.14999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 16018.241
$ws.Range("I99").Value = 19383.4
$ws.Range("K99").Value = 19383.4
$ws.Range("M99").Value = -17885.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4603.1704
$ws.Range("I134").Value = 4633.0454
$ws.Range("K134").Value = 13899.1362
$ws.Range("M134").Value = -11364.1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8335620
$ws.Range("I7").Value = 4372.1665
$ws.Range("J7").Value = 16666868
$ws.Range("K7").Value = 4372.1665
$ws.Range("L7").Value = 16666868
$ws.Range("M7").Value = -4259.1665
$ws.Range("N7").Value = -16667094

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6099.3
$ws.Range("I58").Value = 8501.368
$ws.Range("K58").Value = 8501.368
$ws.Range("M58").Value = -8298.368

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 42138.75
$ws.Range("I107").Value = 82577.75
$ws.Range("K107").Value = 82577.75
$ws.Range("M107").Value = -80657.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1988.8889
$ws.Range("I122").Value = 1650.3334
$ws.Range("K122").Value = 4951.0002
$ws.Range("M122").Value = -2501.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1900023
$ws.Range("I134").Value = 3480611.5
$ws.Range("K134").Value = 10441834.5
$ws.Range("M134").Value = -10439299.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6099.3
$ws.Range("I136").Value = 8501.368
$ws.Range("K136").Value = 25504.104
$ws.Range("M136").Value = -22954.104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1503.35
$ws.Range("I38").Value = 271
$ws.Range("K38").Value = 813
$ws.Range("M38").Value = -466

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 313.25
$ws.Range("I115").Value = 289
$ws.Range("K115").Value = 867
$ws.Range("M115").Value = 308

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4179827.8
$ws.Range("I132").Value = 756.6667
$ws.Range("J132").Value = 4518671.5
$ws.Range("K132").Value = 6810.0003
$ws.Range("L132").Value = 40668043.5
$ws.Range("M132").Value = -4280.0003
$ws.Range("N132").Value = -40673103.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -3900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 7510.515
$ws.Range("J97").Value = 2061.0908
$ws.Range("L97").Value = 2061.0908
$ws.Range("N97").Value = -3053.0908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 21311.691
$ws.Range("I126").Value = 54066
$ws.Range("J126").Value = 11485.4
$ws.Range("K126").Value = 162198
$ws.Range("L126").Value = 34456.2
$ws.Range("M126").Value = -159728
$ws.Range("N126").Value = -39396.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6254.0435
$ws.Range("I132").Value = 5040.4287
$ws.Range("K132").Value = 15121.2861
$ws.Range("M132").Value = -12591.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10848.917
$ws.Range("I22").Value = 14440.154
$ws.Range("J22").Value = 6604.727
$ws.Range("K22").Value = 14440.154
$ws.Range("L22").Value = 6604.727
$ws.Range("M22").Value = -14145.154
$ws.Range("N22").Value = -7194.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 10848.917
$ws.Range("I27").Value = 14440.154
$ws.Range("J27").Value = 6604.727
$ws.Range("K27").Value = 14440.154
$ws.Range("L27").Value = 6604.727
$ws.Range("M27").Value = -14333.154
$ws.Range("N27").Value = -6818.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3506
$ws.Range("I68").Value = 1499.4445
$ws.Range("J68").Value = 6085.857
$ws.Range("K68").Value = 1499.4445
$ws.Range("L68").Value = 6085.857
$ws.Range("M68").Value = -750.4445000000001
$ws.Range("N68").Value = -7583.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3506
$ws.Range("I71").Value = 1499.4445
$ws.Range("J71").Value = 6085.857
$ws.Range("K71").Value = 7497.2225
$ws.Range("L71").Value = 30429.285
$ws.Range("M71").Value = -3753.2225
$ws.Range("N71").Value = -37917.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3459.3
$ws.Range("I93").Value = 3699.7334
$ws.Range("J93").Value = 2738
$ws.Range("K93").Value = 3699.7334
$ws.Range("L93").Value = 2738
$ws.Range("M93").Value = -2451.7334
$ws.Range("N93").Value = -5234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4336.9
$ws.Range("I122").Value = 3833.8333
$ws.Range("K122").Value = 11501.4999
$ws.Range("M122").Value = -9051.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17026.611
$ws.Range("I132").Value = 22259.924
$ws.Range("J132").Value = 3420
$ws.Range("K132").Value = 66779.772
$ws.Range("L132").Value = 10260
$ws.Range("M132").Value = -64249.772
$ws.Range("N132").Value = -15320

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5036.9487
$ws.Range("I136").Value = 1904.8096
$ws.Range("K136").Value = 5714.4288
$ws.Range("M136").Value = -3164.4288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6575.5557
$ws.Range("I122").Value = 2045.875
$ws.Range("K122").Value = 6137.625
$ws.Range("M122").Value = -3687.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 30060.5
$ws.Range("I126").Value = 44200
$ws.Range("K126").Value = 132600
$ws.Range("M126").Value = -130130
